$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.6520907927893234
    "D2" = 0.2482944384878252
    "E2" = 0.2198117960886918
    "F2" = 1.233542482145339
    "G2" = 0.6194693546341625
    "H2" = 0.7520296737587699
    "I2" = 1.028497087521988
    "J2" = 0.2581512450367225
    "K2" = 0.4286591866521974
    "L2" = 0.1905818721582762
    "O2" = 2.728473077746202
    "B3" = 0.6211635282272709
    "D3" = 0.2479319606418215
    "E3" = 0.2213725977103556
    "F3" = 1.241420756502158
    "G3" = 0.6241270415847069
    "H3" = 0.7578060948917127
    "I3" = 1.04073717005058
    "J3" = 0.260297697688222
    "K3" = 0.3742941151593868
    "L3" = 0.1788356043501125
    "O3" = 2.750199244594697
    "B4" = 0.6023144925537451
    "D4" = 0.2477928448671776
    "E4" = 0.2224016038801571
    "F4" = 1.246868195792715
    "G4" = 0.6273740564251398
    "H4" = 0.7616541696155963
    "I4" = 1.048710426779381
    "J4" = 0.2616906733116013
    "K4" = 0.3407701309234596
    "L4" = 0.1716567210223587
    "O4" = 2.764982393132996
    "B5" = 0.5946693865294037
    "D5" = 0.2477572156611743
    "E5" = 0.222838734734335
    "F5" = 1.249241635682559
    "G5" = 0.6287945973221767
    "H5" = 0.7632981393923615
    "I5" = 1.052074759282116
    "J5" = 0.2622772226722092
    "K5" = 0.3270736975882187
    "L5" = 0.168739892989862
    "O5" = 2.771369698015519
    "B6" = 0.5934021196827644
    "D6" = 0.2477525738249469
    "E6" = 0.222912396161604
    "F6" = 1.249645022505696
    "G6" = 0.6290363570362558
    "H6" = 0.763575702607568
    "I6" = 1.052640361292578
    "J6" = 0.2623757611656149
    "K6" = 0.3247973217783056
    "L6" = 0.1682560837326008
    "O6" = 2.772452236928785
    "B7" = 0.6022112411487512
    "D7" = 0.2477922789774212
    "E7" = 0.2224074270545007
    "F7" = 1.246899582874619
    "G7" = 0.6273928201493888
    "H7" = 0.7616760335593611
    "I7" = 1.048755332976883
    "J7" = 0.261698507154512
    "K7" = 0.3405855570107974
    "L7" = 0.1716173484264658
    "O7" = 2.765067064379963
    "B8" = 0.6413983449638749
    "D8" = 0.2481521829893509
    "E8" = 0.220335313003698
    "F8" = 1.236132357099066
    "G8" = 0.6209949646997828
    "H8" = 0.753958888976527
    "I8" = 1.032622482899679
    "J8" = 0.2588757843971425
    "K8" = 0.4099444839504542
    "L8" = 0.1865249686627664
    "O8" = 2.735664813792283
    "B9" = 0.7193309803729449
    "D9" = 0.2495168647651411
    "E9" = 0.2168311633227749
    "F9" = 1.219853501503373
    "G9" = 0.6115213127661718
    "H9" = 0.7412131466291711
    "I9" = 1.004615871092394
    "J9" = 0.2539344891214366
    "K9" = 0.5447811583243833
    "L9" = 0.2160148064047007
    "O9" = 2.68945238869415
    "B10" = 0.7772198421080816
    "D10" = 0.2509172165169105
    "E10" = 0.2145956051178679
    "F10" = 1.210834499383687
    "G10" = 0.6064352264416044
    "H10" = 0.733299869042348
    "I10" = 0.9862482508135724
    "J10" = 0.2506642748290474
    "K10" = 0.643089032554002
    "L10" = 0.2378275139620882
    "O10" = 2.662470788011561
    "B11" = 0.803686011692605
    "D11" = 0.2516398139310354
    "E11" = 0.213651765584844
    "F11" = 1.207368740093116
    "G11" = 0.6045286680850523
    "H11" = 0.7300140469954925
    "I11" = 0.978371102579402
    "J11" = 0.2492543474571168
    "K11" = 0.6876397050553749
    "L11" = 0.2477806022719875
    "O11" = 2.651708446423342
    "B12" = 0.8137264304077974
    "D12" = 0.2519256730413133
    "E12" = 0.2133048408263836
    "F12" = 1.2061478272418
    "G12" = 0.6038652633718726
    "H12" = 0.728814866977217
    "I12" = 0.9754569715399963
    "J12" = 0.2487315884625394
    "K12" = 0.7044845654399978
    "L12" = 0.2515537380197088
    "O12" = 2.647850296510427
    "B13" = 0.8115632468017395
    "D13" = 0.2518635655996917
    "E13" = 0.2133790913354954
    "F13" = 1.206406705152638
    "G13" = 0.6040055342820452
    "H13" = 0.7290711274166597
    "I13" = 0.9760815233027937
    "J13" = 0.2488436784609751
    "K13" = 0.7008578748136358
    "L13" = 0.2507409466075217
    "O13" = 2.648671554085723
    "B14" = 0.8045116814826656
    "D14" = 0.2516630871264312
    "E14" = 0.2136230138699631
    "F14" = 1.207266461790461
    "G14" = 0.6044729154701258
    "H14" = 0.7299144864089371
    "I14" = 0.978129977019357
    "J14" = 0.2492111164396476
    "K14" = 0.6890260592011828
    "L14" = 0.2480909393431148
    "O14" = 2.651386679186544
    "B15" = 0.8001947488686767
    "D15" = 0.2515418783032004
    "E15" = 0.2137737883016229
    "F15" = 1.207804999649682
    "G15" = 0.6047668278362224
    "H15" = 0.730436938385246
    "I15" = 0.979393670019121
    "J15" = 0.2494376341300166
    "K15" = 0.6817753767460033
    "L15" = 0.2464682621116197
    "O15" = 2.653078070499447
    "B16" = 0.7754928115406585
    "D16" = 0.2508717083226912
    "E16" = 0.2146587562917013
    "F16" = 1.211073804262234
    "G16" = 0.6065680189618874
    "H16" = 0.7335209180519726
    "I16" = 0.9867726657520564
    "J16" = 0.2507579785026652
    "K16" = 0.6401740286704865
    "L16" = 0.2371776460446426
    "O16" = 2.663204546256452
    "B17" = 0.7603723000549678
    "D17" = 0.2504824420212799
    "E17" = 0.2152203645766022
    "F17" = 1.213242191250274
    "G17" = 0.6077772833072217
    "H17" = 0.735493212524986
    "I17" = 0.9914219526037833
    "J17" = 0.2515878526914826
    "K17" = 0.6146086367992041
    "L17" = 0.2314857558075403
    "O17" = 2.669803951523647
    "B18" = 0.7516878738032915
    "D18" = 0.2502666052751934
    "E18" = 0.2155502715949265
    "F18" = 1.214549356948829
    "G18" = 0.6085111396547163
    "H18" = 0.7366571797565697
    "I18" = 0.9941411272115062
    "J18" = 0.2520724902233324
    "K18" = 0.599888154444443
    "L18" = 0.2282148054218851
    "O18" = 2.673742049718427
    "B19" = 0.7487496450218032
    "D19" = 0.2501949131609038
    "E19" = 0.2156631557903417
    "F19" = 1.215002243473592
    "G19" = 0.6087661914937499
    "H19" = 0.7370563570579804
    "I19" = 0.9950695270558292
    "J19" = 0.2522378374221208
    "K19" = 0.5949013457017713
    "L19" = 0.2271078191520672
    "O19" = 2.675099863444387
    "B20" = 0.7619806169536787
    "D20" = 0.2505230465775128
    "E20" = 0.21515986807551
    "F20" = 1.213005157197912
    "G20" = 0.6076445889603335
    "H20" = 0.7352801999412719
    "I20" = 0.9909223679891106
    "J20" = 0.2514987541839508
    "K20" = 0.617331774330637
    "L20" = 0.2320913715778801
    "O20" = 2.669086707078165
    "B21" = 0.8065824077760055
    "D21" = 0.2517216412921073
    "E21" = 0.213551083529488
    "F21" = 1.207011448151711
    "G21" = 0.6043340446000016
    "H21" = 0.7296655481810177
    "I21" = 0.9775264301234223
    "J21" = 0.2491028885884727
    "K21" = 0.692502051053367
    "L21" = 0.2488692009735161
    "O21" = 2.650583283812182
    "B22" = 0.8358383069926276
    "D22" = 0.2525762253870028
    "E22" = 0.2125607615525169
    "F22" = 1.20362746433716
    "G22" = 0.6025117990455016
    "H22" = 0.7262588375469505
    "I22" = 0.9691722709735533
    "J22" = 0.24760203079161
    "K22" = 0.7414809788035086
    "L22" = 0.2598583643448364
    "O22" = 2.63975687356799
    "B23" = 0.8202144302092904
    "D23" = 0.2521136256087146
    "E23" = 0.2130837327961466
    "F23" = 1.205384803901453
    "G23" = 0.6034531215312882
    "H23" = 0.7280530372785563
    "I23" = 0.9735943692660634
    "J23" = 0.2483971293603705
    "K23" = 0.7153540025617815
    "L23" = 0.2539911387278835
    "O23" = 2.64541925742688
    "B24" = 0.7612534701715958
    "D24" = 0.2505046644649411
    "E24" = 0.2151871966412617
    "F24" = 1.21311213173616
    "G24" = 0.607704459734407
    "H24" = 0.7353764092633952
    "I24" = 0.991148086166266
    "J24" = 0.2515390121568286
    "K24" = 0.6161007141633661
    "L24" = 0.2318175683299728
    "O24" = 2.669410524700893
    "B25" = 0.698134976763896
    "D25" = 0.2490775518988286
    "E25" = 0.217719463055813
    "F25" = 1.223740381939351
    "G25" = 0.6137551337574649
    "H25" = 0.7444060692869812
    "I25" = 1.011804217279927
    "J25" = 0.2552078498623044
    "K25" = 0.5084342821591861
    "L25" = 0.2080106033999130
    "O25" = 2.700729463117085
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}